$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 5695.4375
$ws.Range("I40").Value = 4796.4
$ws.Range("K40").Value = 4796.4
$ws.Range("M40").Value = -4621.4
$ws.Range("H111").Value = 1182
$ws.Range("I111").Value = 900
$ws.Range("K111").Value = 2700
$ws.Range("M111").Value = 367
$ws.Range("H112").Value = 2411.1177
$ws.Range("J112").Value = 2662.8
$ws.Range("L112").Value = 7988.400000000001
$ws.Range("N112").Value = -10204.4
$ws.Range("H137").Value = 2418.3333
$ws.Range("I137").Value = 1521.7142
$ws.Range("K137").Value = 4565.142599999999
$ws.Range("M137").Value = -2015.142599999999
$ws.Range("H138").Value = 5546.45
$ws.Range("J138").Value = 5812.0527
$ws.Range("L138").Value = 17436.1581
$ws.Range("N138").Value = -27716.1581

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H34").Value = 14999
$ws.Range("J34").Value = 14999
$ws.Range("L34").Value = 14999
$ws.Range("N34").Value = -15541
$ws.Range("H45").Value = 3813.6
$ws.Range("I45").Value = 2688.8
$ws.Range("K45").Value = 2688.8
$ws.Range("M45").Value = -2311.8
$ws.Range("H60").Value = 0
$ws.Range("J60").Value = 0
$ws.Range("L60").Value = 0
$ws.Range("N60").ClearContents()
$ws.Range("H61").Value = 4054.889
$ws.Range("I61").Value = 4054.889
$ws.Range("K61").Value = 4054.889
$ws.Range("M61").Value = -3842.889
$ws.Range("H111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("N111").ClearContents()
$ws.Range("H136").Value = 4054.889
$ws.Range("I136").Value = 4054.889
$ws.Range("K136").Value = 12164.667
$ws.Range("M136").Value = -9614.667000000001

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("N69").ClearContents()
$ws.Range("H72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("N72").ClearContents()
$ws.Range("H99").Value = 1333.7
$ws.Range("I99").Value = 1358.4445
$ws.Range("K99").Value = 1358.4445
$ws.Range("M99").Value = 139.5554999999999
$ws.Range("H107").Value = 4642.6313
$ws.Range("I107").Value = 1518.3334
$ws.Range("K107").Value = 1518.3334
$ws.Range("M107").Value = 401.6666

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H19").Value = 209.53847
$ws.Range("I19").Value = 130.44444
$ws.Range("K19").Value = 130.44444
$ws.Range("M19").Value = 39.55556000000001
$ws.Range("H24").Value = 209.53847
$ws.Range("I24").Value = 130.44444
$ws.Range("K24").Value = 130.44444
$ws.Range("M24").Value = 39.55556000000001
$ws.Range("H57").Value = 10000
$ws.Range("I57").Value = 10000
$ws.Range("K57").Value = 10000
$ws.Range("M57").Value = -9440
$ws.Range("H86").Value = 2998.75
$ws.Range("J86").Value = 3000
$ws.Range("L86").Value = 3000
$ws.Range("N86").Value = -5246
$ws.Range("H89").Value = 2998.75
$ws.Range("J89").Value = 3000
$ws.Range("L89").Value = 15000
$ws.Range("N89").Value = -26232
$ws.Range("H93").Value = 9475.833000000001
$ws.Range("I93").Value = 5281.4
$ws.Range("K93").Value = 5281.4
$ws.Range("M93").Value = -3409.4
$ws.Range("H95").Value = 10000
$ws.Range("J95").Value = 10000
$ws.Range("L95").Value = 10000
$ws.Range("N95").Value = -15492
$ws.Range("H96").Value = 16000
$ws.Range("J96").Value = 16000
$ws.Range("L96").Value = 16000
$ws.Range("N96").Value = -21492

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 8170.923
$ws.Range("J39").Value = 8170.923
$ws.Range("L39").Value = 24512.769
$ws.Range("N39").Value = -25100.769
$ws.Range("H40").Value = 205.18182
$ws.Range("J40").Value = 421.8
$ws.Range("L40").Value = 1687.2
$ws.Range("N40").Value = -1825.2
$ws.Range("H141").Value = 1152.5
$ws.Range("I141").Value = 1152.5
$ws.Range("K141").Value = 3457.5
$ws.Range("M141").Value = 1722.5

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H7").Value = 30545636
$ws.Range("I7").Value = 20000250
$ws.Range("J7").Value = 58666668
$ws.Range("K7").Value = 20000250
$ws.Range("L7").Value = 58666668
$ws.Range("M7").Value = -20000138
$ws.Range("N7").Value = -58666892
$ws.Range("H8").Value = 30545636
$ws.Range("I8").Value = 20000250
$ws.Range("J8").Value = 58666668
$ws.Range("K8").Value = 20000250
$ws.Range("L8").Value = 58666668
$ws.Range("M8").Value = -20000111
$ws.Range("N8").Value = -58666946
$ws.Range("H10").Value = 9333
$ws.Range("I10").Value = 9333
$ws.Range("K10").Value = 9333
$ws.Range("M10").Value = -9164
$ws.Range("H12").Value = 999.5
$ws.Range("I12").Value = 999.5
$ws.Range("K12").Value = 999.5
$ws.Range("M12").Value = -859.5
$ws.Range("H14").Value = 751.5
$ws.Range("I14").Value = 751.5
$ws.Range("K14").Value = 751.5
$ws.Range("M14").Value = -583.5

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H14").Value = 5000
$ws.Range("J14").Value = 5000
$ws.Range("L14").Value = 5000
$ws.Range("N14").Value = -5344
$ws.Range("H16").Value = 2461.7273
$ws.Range("J16").Value = 2248.5
$ws.Range("L16").Value = 2248.5
$ws.Range("N16").Value = -2588.5
$ws.Range("H40").Value = 6078.5
$ws.Range("I40").Value = 5598.125
$ws.Range("K40").Value = 5598.125
$ws.Range("M40").Value = -5462.125
$ws.Range("H45").Value = 0
$ws.Range("I45").Value = 0
$ws.Range("K45").Value = 0
$ws.Range("M45").ClearContents()
$ws.Range("H46").Value = 9422.333000000001
$ws.Range("I46").Value = 5000
$ws.Range("J46").Value = 9975.125
$ws.Range("K46").Value = 5000
$ws.Range("L46").Value = 9975.125
$ws.Range("M46").Value = -4812
$ws.Range("N46").Value = -10351.125
$ws.Range("H61").Value = 5828.9
$ws.Range("I61").Value = 4464.8335
$ws.Range("K61").Value = 4464.8335
$ws.Range("M61").Value = -4262.8335
$ws.Range("H100").Value = 8749.75
$ws.Range("I100").Value = 5000
$ws.Range("K100").Value = 5000
$ws.Range("M100").Value = -4459
$ws.Range("H113").Value = 5828.9
$ws.Range("I113").Value = 4464.8335
$ws.Range("K113").Value = 4464.8335
$ws.Range("M113").Value = -2294.8335
$ws.Range("H132").Value = 3678.7144
$ws.Range("I132").Value = 3562.75
$ws.Range("K132").Value = 10688.25
$ws.Range("M132").Value = -8158.25
$ws.Range("H136").Value = 1902
$ws.Range("I136").Value = 1902
$ws.Range("K136").Value = 5706
$ws.Range("M136").Value = -3156

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H8").Value = 0
$ws.Range("J8").Value = 0
$ws.Range("L8").Value = 0
$ws.Range("N8").ClearContents()
$ws.Range("H11").Value = 2999.6667
$ws.Range("I11").Value = 1999
$ws.Range("K11").Value = 1999
$ws.Range("M11").Value = -1857
$ws.Range("H12").Value = 9000
$ws.Range("J12").Value = 9000
$ws.Range("L12").Value = 9000
$ws.Range("N12").Value = -9284
$ws.Range("H29").Value = 4555.5557
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 4555.5557
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = 4555.5557
$ws.Range("N29").Value = -5135.5557
$ws.Range("M29").ClearContents()
$ws.Range("H54").Value = 6125
$ws.Range("J54").Value = 6125
$ws.Range("L54").Value = 6125
$ws.Range("N54").Value = -7165
$ws.Range("H114").Value = 30000
$ws.Range("J114").Value = 30000
$ws.Range("L114").Value = 30000
$ws.Range("N114").Value = -38678
$ws.Range("H136").Value = 3028.9092
$ws.Range("I136").Value = 2365.3333
$ws.Range("J136").Value = 3825.2
$ws.Range("K136").Value = 7095.999899999999
$ws.Range("L136").Value = 11475.6
$ws.Range("M136").Value = -4545.999899999999
$ws.Range("N136").Value = -16575.6
